$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.120.19"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.820.20"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "311.75"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.4456"
$ws.Range("E7").Value = "  +5.19%  "
$ws.Range("D8").Value = "0.3739"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "0.07439"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "0.8719"
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("D11").Value = "20.88"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.809.48"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "6.736"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "94.32"
$ws.Range("E14").Value = "  +5.30%  "
$ws.Range("D15").Value = "5.332"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "0.07082"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "0.000008749"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "0.9989"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "15.00"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "27.141.87"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "5.216"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").Value = "2.036.94"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "2.369"
$ws.Range("E26").Value = "  +5.00%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "151.54"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "5.315"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "118.04"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08817"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7652"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.177"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.568"
$ws.Range("E34").Value = "  +3.29%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.885"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "0.9984"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.099"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01984"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05272"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.376"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5306"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").Value = "0.1719"
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.852"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.170"
$ws.Range("E44").Value = "  +10.19%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.699"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5042"
$ws.Range("E46").Value = "  +6.75%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.63"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.707"
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "105.65"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "0.9981"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06365"
$ws.Range("E51").Value = "  +0.63%  "
